$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 209, shifting the existing 209:217 block
# (weekly Coliflor price records) down to 211:219.
$ws.Rows("209:210").Insert()

# Row 209 - new weekly record (Primera quality), date 2021-11-09 (serial 44509)
$ws.Range("A209").Value = 4
$ws.Range("B209").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C209").Value = "Los Lagos"
$ws.Range("D209").Value = 44509
$ws.Range("E209").Value = 10
$ws.Range("F209").Value = 100112008
$ws.Range("G209").Value = "Coliflor"
$ws.Range("H209").Value = "Sin especificar"
$ws.Range("I209").Value = "Primera"
$ws.Range("J209").Value = 600
$ws.Range("K209").Value = 1200
$ws.Range("L209").Value = 1200
$ws.Range("M209").Value = 1200
$ws.Range("N209").Value = "$/unidad"
$ws.Range("O209").Value = "Región Metropolitana"
$ws.Range("P209").Value = 1200
$ws.Range("Q209").Value = 1
$ws.Range("R209").Value = "Hortaliza"

# Row 210 - new weekly record (Segunda quality), same date
$ws.Range("A210").Value = 4
$ws.Range("B210").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C210").Value = "Los Lagos"
$ws.Range("D210").Value = 44509
$ws.Range("E210").Value = 10
$ws.Range("F210").Value = 100112008
$ws.Range("G210").Value = "Coliflor"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Segunda"
$ws.Range("J210").Value = 600
$ws.Range("K210").Value = 1000
$ws.Range("L210").Value = 1000
$ws.Range("M210").Value = 1000
$ws.Range("N210").Value = "$/unidad"
$ws.Range("O210").Value = "Región Metropolitana"
$ws.Range("P210").Value = 1000
$ws.Range("Q210").Value = 1
$ws.Range("R210").Value = "Hortaliza"
